# Weekly data refresh: insert a new row of "Arveja Verde" price data at the
# top of the historical block (row 77), pushing the existing rows 77-98
# down to 78-99.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(77).Insert()

$ws.Cells.Item(77, 1).Value  = 7
$ws.Cells.Item(77, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(77, 3).Value  = "Ñuble"
$ws.Cells.Item(77, 4).Value  = 44932
$ws.Cells.Item(77, 5).Value  = 16
$ws.Cells.Item(77, 6).Value  = 100112022
$ws.Cells.Item(77, 7).Value  = "Arveja Verde"
$ws.Cells.Item(77, 8).Value  = "Sin especificar"
$ws.Cells.Item(77, 9).Value  = "Primera"
$ws.Cells.Item(77, 10).Value = 60
$ws.Cells.Item(77, 11).Value = 20000
$ws.Cells.Item(77, 12).Value = 21000
$ws.Cells.Item(77, 13).Value = 20500
$ws.Cells.Item(77, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(77, 15).Value = "Región de Ñuble"
$ws.Cells.Item(77, 16).Value = 820
$ws.Cells.Item(77, 17).Value = 25
$ws.Cells.Item(77, 18).Value = "Hortaliza"
